# Attorneys.xlsx: insert a new header-ish row ("an" / "attorney") right
# below the title row, pushing the existing attorney list down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 (existing rows 2..24 become 3..25).
$ws.Rows("2:2").Insert()

# New row 2 content.
$ws.Range("A2").Value = "an"
$ws.Range("B2").Value = "attorney"

# The inserted row picked up row 1's (bold, header) formatting via the
# Insert() copy-down behaviour; restore plain (non-bold) formatting and a
# text number format on B2, matching the rest of the data rows.
$ws.Range("A2:B2").Font.Bold = $false
$ws.Range("B2").NumberFormat = "@"

# Column A width matches the rest of the (now shifted) name column.
$ws.Columns("A").ColumnWidth = 17.17

# Leave the selection on B3 (the first real data row, directly below the
# new row), matching the saved file's cursor position.
[void]$ws.Range("B3").Select()
